# Excel Skills for Business Intermediate I - Week 2 - Changing Case
# Update the "Full Name" (D) formulas to PROPER-case the concatenated
# first/last names, and the "Email" (E) formulas to lower-case the
# generated email address, for every employee row (4-38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 4; $r -le 38; $r++) {
    $ws.Range("D$r").Formula = "=PROPER(CONCATENATE(C$r, "" "", B$r))"
    $ws.Range("E$r").Formula = "=LOWER(C$r&"".""&B$r&""@pushpin.com"")"
}

# Match the cursor position left behind in the saved workbook.
$ws.Range("O6:O7").Select() | Out-Null
